$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (vNameAccount), shifting C:G to D:H
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted "Status" column
$ws.Range("C1").Value = "Status"

# Status values for each data row
$ws.Range("C2").Value = "Failed"
$ws.Range("C3").Value = "Failed"
$ws.Range("C4").Value = "Failed"
$ws.Range("C5").Value = ""
$ws.Range("C6").Value = ""

# RunTest column (B) changes: rows 3 and 4 flip from Yes to No
$ws.Range("B3").Value = "No"
$ws.Range("B4").Value = "No"

# vOutData column (now H) timestamps updated for rows 2-4
$ws.Range("H2").Value = "25_04_2020--19_26_48 617"
$ws.Range("H3").Value = "25_04_2020--19_14_17 408"
$ws.Range("H4").Value = "25_04_2020--19_14_48 020"

# Column widths (Status column matches RunTest's width; new vOutData timestamp
# column is widened to fit the longer timestamp strings)
$ws.Columns.Item(3).ColumnWidth = 9
$ws.Columns.Item(8).ColumnWidth = 24.5

# Selection matches the final state
$ws.Range("B5").Select()
